$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on numeric-looking Price cells so values are preserved as strings
$textCells = @("D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D21", "D23", "D24", "D25", "D26", "D27", "D29", "D30", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.875.15"
$ws.Range("E2").Value = "  -0.94%  "

$ws.Range("D3").Value = "1.860.00"
$ws.Range("E3").Value = "  -0.44%  "

$ws.Range("D4").Value = "0.9995"
$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").Value = "304.52"
$ws.Range("E5").Value = "  -0.83%  "

$ws.Range("D6").Value = "0.9995"
$ws.Range("E6").Value = "  -0.13%  "

$ws.Range("D7").Value = "0.5051"
$ws.Range("E7").Value = "  -1.42%  "

$ws.Range("D8").Value = "0.3628"
$ws.Range("E8").Value = "  -2.90%  "

$ws.Range("D9").Value = "0.07180"
$ws.Range("E9").Value = "  +0.60%  "

$ws.Range("D10").Value = "0.8937"
$ws.Range("E10").Value = "  +0.67%  "

$ws.Range("D11").Value = "20.72"
$ws.Range("E11").Value = "  +0.37%  "

$ws.Range("D12").Value = "1.867.76"
$ws.Range("E12").Value = "  -2.40%  "

$ws.Range("D13").Value = "0.07498"
$ws.Range("E13").Value = "  -0.52%  "

$ws.Range("D14").Value = "92.16"
$ws.Range("E14").Value = "  +3.37%  "

$ws.Range("D15").Value = "5.224"
$ws.Range("E15").Value = "  -1.70%  "

$ws.Range("D16").Value = "0.9998"
$ws.Range("E16").Value = "  -0.10%  "

$ws.Range("D17").Value = "0.000008479"
$ws.Range("E17").Value = "  +0.09%  "

$ws.Range("D18").Value = "14.15"
$ws.Range("E18").Value = "  +0.27%  "

$ws.Range("D19").Value = "0.9994"
$ws.Range("E19").Value = "  -0.16%  "

$ws.Range("D20").Value = "26.914.39"
$ws.Range("E20").Value = "  -1.00%  "

$ws.Range("D21").Value = "5.032"
$ws.Range("E21").Value = "  -0.28%  "

$ws.Range("D22").Value = "2.107.35"
$ws.Range("E22").Value = "  +0.73%  "

$ws.Range("D23").Value = "10.35"
$ws.Range("E23").Value = "  -1.96%  "

$ws.Range("D24").Value = "6.395"
$ws.Range("E24").Value = "  -1.12%  "

$ws.Range("D25").Value = "147.61"
$ws.Range("E25").Value = "  -1.70%  "

$ws.Range("D26").Value = "1.791"
$ws.Range("E26").Value = "  -2.54%  "

$ws.Range("D27").Value = "17.86"
$ws.Range("E27").Value = "  -0.31%  "

$ws.Range("E28").Value = "  -1.24%  "

$ws.Range("D29").Value = "113.07"
$ws.Range("E29").Value = "  +0.36%  "

$ws.Range("D30").Value = "4.687"
$ws.Range("E30").Value = "  -1.26%  "

$ws.Range("E31").Value = "  +0.20%  "

$ws.Range("D32").Value = "0.09256"
$ws.Range("E32").Value = "  +2.72%  "

$ws.Range("D33").Value = "0.05094"
$ws.Range("E33").Value = "  -0.65%  "

$ws.Range("D34").Value = "0.7481"
$ws.Range("E34").Value = "  +1.96%  "

$ws.Range("D35").Value = "2.979"
$ws.Range("E35").Value = "  -3.90%  "

$ws.Range("D36").Value = "1.148"
$ws.Range("E36").Value = "  -0.95%  "

$ws.Range("D37").Value = "3.284"
$ws.Range("E37").Value = "  +7.60%  "

$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "2.515"
$ws.Range("E38").Value = "  +0.68%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "0.02001"
$ws.Range("E39").Value = "  -2.13%  "

$ws.Range("D40").Value = "0.5495"
$ws.Range("E40").Value = "  +3.57%  "

$ws.Range("E41").Value = "  -0.86%  "

$ws.Range("D42").Value = "118.06"
$ws.Range("E42").Value = "  +1.41%  "

$ws.Range("D43").Value = "6.498"
$ws.Range("E43").Value = "  -1.49%  "

$ws.Range("D44").Value = "8.505"
$ws.Range("E44").Value = "  +2.01%  "

$ws.Range("D45").Value = "0.1468"
$ws.Range("E45").Value = "  -0.31%  "

$ws.Range("D46").Value = "0.4667"
$ws.Range("E46").Value = "  +1.09%  "

$ws.Range("D47").Value = "0.9990"
$ws.Range("E47").Value = "  -0.15%  "

$ws.Range("D48").Value = "10.08"
$ws.Range("E48").Value = "  +1.17%  "

$ws.Range("D49").Value = "1.564"
$ws.Range("E49").Value = "  -0.25%  "

$ws.Range("D50").Value = "36.87"
$ws.Range("E50").Value = "  +1.16%  "

$ws.Range("D51").Value = "63.18"
$ws.Range("E51").Value = "  -2.10%  "
